$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Caso1")
$ws.Range("G2").Value = 0.9971016645431519
$ws.Range("G3").Value = 0.997994065284729
$ws.Range("G4").Value = 0.9983851313591003
$ws.Range("G5").Value = 0.9982051849365234
$ws.Range("G6").Value = 0.9981799125671387
$ws.Range("G7").Value = 0.9986467361450195
$ws.Range("G8").Value = 0.9984886646270752
$ws.Range("G9").Value = 0.9972221255302429
$ws.Range("G10").Value = 0.9993353486061096
$ws.Range("G11").Value = 0.9978989362716675
$ws.Range("G12").Value = 0.9978556036949158
$ws.Range("G13").Value = 0.9979280233383179
$ws.Range("G14").Value = 0.9994798302650452
$ws.Range("G15").Value = 0.999695360660553
$ws.Range("G16").Value = 0.9974928498268127
$ws.Range("G17").Value = 0.998495876789093
$ws.Range("G18").Value = 0.9980384111404419
$ws.Range("G19").Value = 0.997744619846344

$ws = $wb.Worksheets.Item("Caso2")
$ws.Range("G2").Value = 0.9957005977630615
$ws.Range("G3").Value = 0.9965415000915527
$ws.Range("G4").Value = 0.9969353675842285
$ws.Range("G5").Value = 0.9967901110649109
$ws.Range("G6").Value = 0.9967455267906189
$ws.Range("G7").Value = 0.9972108602523804
$ws.Range("G8").Value = 0.9970960021018982
$ws.Range("G9").Value = 0.9958401918411255
$ws.Range("G10").Value = 0.9978943467140198
$ws.Range("G11").Value = 0.9964532852172852
$ws.Range("G12").Value = 0.9964525103569031
$ws.Range("G13").Value = 0.996525764465332
$ws.Range("G14").Value = 0.9981326460838318
$ws.Range("G15").Value = 0.9983751177787781
$ws.Range("G16").Value = 0.9960439205169678
$ws.Range("G17").Value = 0.9970303177833557
$ws.Range("G18").Value = 0.9965879321098328
$ws.Range("G19").Value = 0.9963411092758179

$ws = $wb.Worksheets.Item("Caso3")
$ws.Range("G2").Value = 0.9940349459648132
$ws.Range("G3").Value = 0.9948811531066895
$ws.Range("G4").Value = 0.9952586889266968
$ws.Range("G5").Value = 0.9950906038284302
$ws.Range("G6").Value = 0.9950611591339111
$ws.Range("G7").Value = 0.9955226182937622
$ws.Range("G8").Value = 0.9953632950782776
$ws.Range("G9").Value = 0.9941365718841553
$ws.Range("G10").Value = 0.9961778521537781
$ws.Range("G11").Value = 0.9948064684867859
$ws.Range("G12").Value = 0.9947551488876343
$ws.Range("G13").Value = 0.9948323369026184
$ws.Range("G14").Value = 0.9964390397071838
$ws.Range("G15").Value = 0.9966772794723511
$ws.Range("G16").Value = 0.9943928122520447
$ws.Range("G17").Value = 0.9953623414039612
$ws.Range("G18").Value = 0.9949227571487427
$ws.Range("G19").Value = 0.9946147799491882

